# Actualizacion Datos Personales 4 nov
$wb = $excel.ActiveWorkbook

# ---- Sheet "3ASV" ----
$ws1 = $wb.Worksheets.Item("3ASV")

# Row 8 (GARCIA LEON, JESUS SAMUEL): correct student e-mail
$ws1.Range("E8").Value = "erikaleonpalacios7@gmail.com"

# Row 4 (Tutor name): fix typo ITSEL -> ITZEL
$ws1.Range("H4").Value = "GUADALUPE ITZEL TEPEPA ROSAS"

# Row 15 (Tutor name): remove duplicated "MOLINA MORALES"
$ws1.Range("H15").Value = "MIRIAM MOLINA MORALES"

# Row 19 (RIVERA VARGAS, KELLY ITZEL): add missing tutor e-mail
$ws1.Range("I19").Value = "Kelly_rivera_vargas@gmail.com"

# ---- Sheet "5AEV" ----
$ws2 = $wb.Worksheets.Item("5AEV")

# Row 21 (ROMERO CORTES, ARTURO): fix typo gamil.com -> gmail.com
$ws2.Range("E21").Value = "arturombappe10@gmail.com"

# Row 6: remove the erroneous tutor e-mail "Marielc"
$ws2.Range("I6").ClearContents()

# Row 23 (SOLANO NOGALES, AGUSTIN): add missing tutor name
$ws2.Range("H23").Value = "AGUSTÍN SOLANO MARTÍNEZ"
